$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values that changed
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "247.06"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.75"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.292"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05723"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.441"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8101"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8725"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1429"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07380"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03051"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03128"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09393"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.954"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001570"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04811"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005854"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006148"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005123"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009961"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0001500"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.744"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.315"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3277"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03928"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006703"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1068"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002550"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007479"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005629"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6004"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1757"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01011"

# Update Hora (column G) values from 5 to 6 for all data rows (2-51)
for ($row = 2; $row -le 51; $row++) {
    $ws.Range("G$row").NumberFormat = "@"
    $ws.Range("G$row").Value = "6"
}